# ADONetEssentials.pptx update
#  1. Bump the cached "datetimeFigureOut" date placeholder text on the
#     slide master and every slide layout from 26-03-2023 to 01-04-2024.
#  2. Bump the Microsoft.Data.SqlClient documentation link on the last
#     slide from the v5.1 docs URL to the v5.2 docs URL (collapsing the
#     four runs that made up the link into a single run).

$p = $ppt.ActivePresentation

$oldDate = "26-03-2023"
$newDate = "01-04-2024"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master date placeholder.
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every custom layout's date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# Last slide: update the Microsoft.Data.SqlClient hyperlink text,
# merging the split runs into one run with the new version number.
$lastSlide = $p.Slides.Item($p.Slides.Count)
for ($i = 1; $i -le $lastSlide.Shapes.Count; $i++) {
    $shp = $lastSlide.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -like "*sqlclient-dotnet-core-5.1*") {
            $oldUrl = "https://learn.microsoft.com/en-us/dotnet/api/microsoft.data.sqlclient?view=sqlclient-dotnet-core-5.1"
            $newUrl = "https://learn.microsoft.com/en-us/dotnet/api/microsoft.data.sqlclient?view=sqlclient-dotnet-core-5.2"
            $sub = $tr.Characters(1, $oldUrl.Length)
            $sub.Text = $newUrl
        }
    }
}
